$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# HOUR_APPR_PROCESS_START (column V) values, converted from plain hour
# numbers into "HH:00:00" formatted time strings (stored as text, not as
# a numeric/time value).
$hours = @{
    2  = 13
    3  = 13
    4  = 10
    5  = 13
    6  = 10
    7  = 15
    8  = 14
    9  = 11
    10 = 11
    11 = 13
    12 = 12
}

foreach ($row in $hours.Keys) {
    $cell = $ws.Range("V$row")
    $cell.Value = "{0}:00:00" -f $hours[$row]
}
